$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cell values to match the re-shuffled weekly data dump ---
$ws.Range("D2").Value = 44350
$ws.Range("D3").Value = 44350
$ws.Range("D4").Value = 44349
$ws.Range("O4").Value = "Región Metropolitana"
$ws.Range("D5").Value = 44349
$ws.Range("O5").Value = "Región Metropolitana"
$ws.Range("D6").Value = 44327
$ws.Range("D7").Value = 44327
$ws.Range("D8").Value = 44280
$ws.Range("D9").Value = 44280
$ws.Range("D10").Value = 44383
$ws.Range("D11").Value = 44383
$ws.Range("D12").Value = 44336
$ws.Range("D13").Value = 44336
$ws.Range("D14").Value = 44306
$ws.Range("D15").Value = 44306
$ws.Range("D16").Value = 44160
$ws.Range("D17").Value = 44160
$ws.Range("D18").Value = 44285
$ws.Range("D19").Value = 44285
$ws.Range("D20").Value = 44425
$ws.Range("D21").Value = 44425
$ws.Range("D22").Value = 44386
$ws.Range("D23").Value = 44386
$ws.Range("D24").Value = 44166
$ws.Range("D25").Value = 44166
$ws.Range("D26").Value = 44405
$ws.Range("D27").Value = 44405
$ws.Range("D28").Value = 44308
$ws.Range("D29").Value = 44308
$ws.Range("D30").Value = 44231
$ws.Range("D31").Value = 44231
$ws.Range("D32").Value = 44355
$ws.Range("O32").Value = "Región de Ñuble"
$ws.Range("D33").Value = 44355
$ws.Range("O33").Value = "Región de Ñuble"
$ws.Range("D34").Value = 44203
$ws.Range("D35").Value = 44203
$ws.Range("D36").Value = 44239
$ws.Range("D37").Value = 44239
$ws.Range("D38").Value = 44330
$ws.Range("D39").Value = 44330
$ws.Range("D40").Value = 44187
$ws.Range("D41").Value = 44187
$ws.Range("D42").Value = 44344
$ws.Range("N42").Value = "`$/docena de 1 kilo"
$ws.Range("D43").Value = 44344
$ws.Range("N43").Value = "`$/docena de 1 kilo"
$ws.Range("D44").Value = 44320
$ws.Range("D45").Value = 44320
$ws.Range("D46").Value = 44278
$ws.Range("J46").Value = 300
$ws.Range("D47").Value = 44278
$ws.Range("J47").Value = 150
$ws.Range("D48").Value = 44250
$ws.Range("O48").Value = "Región de Arica y Parinacota"
$ws.Range("D49").Value = 44250
$ws.Range("O49").Value = "Región de Arica y Parinacota"
$ws.Range("D50").Value = 44334
$ws.Range("D51").Value = 44334
$ws.Range("D52").Value = 44299
$ws.Range("D53").Value = 44299
$ws.Range("D54").Value = 44252
$ws.Range("D55").Value = 44252
$ws.Range("D56").Value = 44292
$ws.Range("D57").Value = 44292
$ws.Range("D58").Value = 44224
$ws.Range("D59").Value = 44224
$ws.Range("D60").Value = 44217
$ws.Range("D61").Value = 44217
$ws.Range("D62").Value = 44362
$ws.Range("D63").Value = 44362
$ws.Range("D64").Value = 44168
$ws.Range("D65").Value = 44168
$ws.Range("D66").Value = 44272
$ws.Range("D67").Value = 44272
$ws.Range("D68").Value = 44365
$ws.Range("J68").Value = 200
$ws.Range("D69").Value = 44365
$ws.Range("J69").Value = 100
$ws.Range("D72").Value = 44433
$ws.Range("O72").Value = "Región de Ñuble"
$ws.Range("D73").Value = 44433
$ws.Range("O73").Value = "Región de Ñuble"
$ws.Range("D74").Value = 44274
$ws.Range("D75").Value = 44274
$ws.Range("D76").Value = 44194
$ws.Range("D77").Value = 44194
$ws.Range("D78").Value = 44222
$ws.Range("N78").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("D79").Value = 44222
$ws.Range("N79").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("D80").Value = 44398
$ws.Range("D81").Value = 44398
$ws.Range("D82").Value = 44316
$ws.Range("D83").Value = 44316
$ws.Range("D84").Value = 44341
$ws.Range("D85").Value = 44341
$ws.Range("D86").Value = 44237
$ws.Range("D87").Value = 44237
$ws.Range("D88").Value = 44391
$ws.Range("D89").Value = 44391
$ws.Range("D90").Value = 44313
$ws.Range("D91").Value = 44313
$ws.Range("D92").Value = 44442
$ws.Range("J92").Value = 300
$ws.Range("D93").Value = 44442
$ws.Range("J93").Value = 150
$ws.Range("D94").Value = 44435
$ws.Range("D95").Value = 44435

# --- Append two new rows (96, 97) for the latest "Perejil" price report ---
$ws.Range("A96").Value = 11
$ws.Range("B96").Value = "Vega Monumental Concepción"
$ws.Range("C96").Value = "Bíobío"
$ws.Range("D96").Value = 44400
$ws.Range("E96").Value = 8
$ws.Range("F96").Value = 100112044
$ws.Range("G96").Value = "Perejil"
$ws.Range("H96").Value = "Sin especificar"
$ws.Range("I96").Value = "Primera"
$ws.Range("J96").Value = 200
$ws.Range("K96").Value = 600
$ws.Range("L96").Value = 700
$ws.Range("M96").Value = 650
$ws.Range("N96").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O96").Value = "Región de Ñuble"
$ws.Range("P96").Value = 650
$ws.Range("Q96").Value = 1
$ws.Range("R96").Value = "Hortaliza"
$ws.Range("D96").NumberFormat = $ws.Range("D2").NumberFormat

$ws.Range("A97").Value = 11
$ws.Range("B97").Value = "Vega Monumental Concepción"
$ws.Range("C97").Value = "Bíobío"
$ws.Range("D97").Value = 44400
$ws.Range("E97").Value = 8
$ws.Range("F97").Value = 100112044
$ws.Range("G97").Value = "Perejil"
$ws.Range("H97").Value = "Sin especificar"
$ws.Range("I97").Value = "Segunda"
$ws.Range("J97").Value = 100
$ws.Range("K97").Value = 500
$ws.Range("L97").Value = 500
$ws.Range("M97").Value = 500
$ws.Range("N97").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O97").Value = "Región de Ñuble"
$ws.Range("P97").Value = 500
$ws.Range("Q97").Value = 1
$ws.Range("R97").Value = "Hortaliza"
$ws.Range("D97").NumberFormat = $ws.Range("D2").NumberFormat

